# Apply updates to the stimuli worksheet:
#  - practice rows (2-5): fill in the "carrier" column D to mirror column K
#  - generic rows (6-9): fill in the "pair_kind" column J with the unique
#    video/audio label
#  - new unique rows (14-21): fill in "kind" column C and "carrier" column D

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "can"
    "D3" = "where"
    "D4" = "do"
    "D5" = "look"

    "J6" = "unique_video"
    "J7" = "unique_video"
    "J8" = "unique_audio"
    "J9" = "unique_audio"

    "C14" = "unique_video"
    "D14" = "look"
    "C15" = "unique_video"
    "D15" = "look"
    "C16" = "unique_video"
    "D16" = "where"
    "C17" = "unique_video"
    "D17" = "where"

    "C18" = "unique_audio"
    "D18" = "can"
    "C19" = "unique_audio"
    "D19" = "can"
    "C20" = "unique_audio"
    "D20" = "do"
    "C21" = "unique_audio"
    "D21" = "do"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
